$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last header cell (E1) into the new header cell F1,
# then set its text, so F1 gets the same style (bold, border, centered) as
# the other header cells without creating a new style entry.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update the existing numeric values in row 2
$ws.Range("B2").Value = 0.0557715238619978
$ws.Range("C2").Value = 0.9992448879935781
$ws.Range("D2").Value = 0.1901453129850954

# Add the new value cell F2 (model description)
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"
